$wb = $excel.ActiveWorkbook

# --- Rename "Items" -> "Actual" ---
$actual = $wb.Worksheets.Item(1)
$actual.Name = "Actual"

# --- Insert a new "Estimated" sheet between Actual and Analysis ---
$estimated = $wb.Worksheets.Add($null, $actual)
$estimated.Name = "Estimated"

# Fill in the Estimated sheet content (order chosen to match the
# original authoring order of the shared-string table).
$estimated.Range("A5").Value = "Other prototype parts"
$estimated.Range("B5").Value = 500
$estimated.Range("A4").Value = "Manufacturability Consultation (2) (estimate)"
$estimated.Range("B4").Value = 400
$estimated.Range("A2").Value = "PCB Design"
$estimated.Range("B2").Value = 400
$estimated.Range("A3").Value = "PCB Fabrication (10) (estimate)"
$estimated.Range("B3").Value = 250
$estimated.Range("A6").Value = "Small batch manufacturing of device"
$estimated.Range("A8").Value = "Total"
$estimated.Range("B8").Formula = "=SUM(B2:B7)"
$estimated.Range("A7").Value = "Factor of Safety"
$estimated.Range("B7").Value = 1000
$estimated.Range("B1").Value = "Cost ($)"

# Column A width (best-fit-ish, matches the source's ~42-wide column)
$estimated.Columns.Item(1).ColumnWidth = 41

# Match the saved selection on the Estimated sheet
$estimated.Range("A9").Select() | Out-Null
